# Scheduled-runner refresh of the Leve profit calculations across all job sheets.
# Mirrors an upstream Universalis price pull: updates currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ) and the derived LeveProfit(NQ/HQ) columns (H:N) for the rows
# whose market data changed since the last run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1424.625
$ws.Range("I40").Value = 1399.9
$ws.Range("J40").Value = 1465.8334
$ws.Range("K40").Value = 1399.9
$ws.Range("L40").Value = 1465.8334
$ws.Range("M40").Value = -1224.9
$ws.Range("N40").Value = -1815.8334
$ws.Range("H132").Value = 1399.2153
$ws.Range("I132").Value = 1309.3833
$ws.Range("K132").Value = 3928.1499
$ws.Range("M132").Value = -1398.1499
$ws.Range("H137").Value = 1475.8334
$ws.Range("I137").Value = 1479.8422
$ws.Range("J137").Value = 1460.6
$ws.Range("K137").Value = 4439.5266
$ws.Range("L137").Value = 4381.799999999999
$ws.Range("M137").Value = -1889.5266
$ws.Range("N137").Value = -9481.8
$ws.Range("H138").Value = 2337.7317
$ws.Range("I138").Value = 1302.4
$ws.Range("J138").Value = 3596.919
$ws.Range("K138").Value = 3907.2
$ws.Range("L138").Value = 10790.757
$ws.Range("M138").Value = 1232.8
$ws.Range("N138").Value = -21070.757

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3312.625
$ws.Range("I26").Value = 3214.4285
$ws.Range("J26").Value = 4000
$ws.Range("K26").Value = 3214.4285
$ws.Range("L26").Value = 4000
$ws.Range("M26").Value = -2884.4285
$ws.Range("N26").Value = -4660
$ws.Range("H32").Value = 11230.179
$ws.Range("I32").Value = 11008.471
$ws.Range("K32").Value = 11008.471
$ws.Range("M32").Value = -10721.471
$ws.Range("H45").Value = 1670.909
$ws.Range("I45").Value = 1542.8572
$ws.Range("J45").Value = 1895
$ws.Range("K45").Value = 1542.8572
$ws.Range("L45").Value = 1895
$ws.Range("M45").Value = -1165.8572
$ws.Range("N45").Value = -2649
$ws.Range("H61").Value = 1831.6552
$ws.Range("I61").Value = 1751.9445
$ws.Range("J61").Value = 1962.091
$ws.Range("K61").Value = 1751.9445
$ws.Range("L61").Value = 1962.091
$ws.Range("M61").Value = -1539.9445
$ws.Range("N61").Value = -2386.091
$ws.Range("H74").Value = 809.2727
$ws.Range("I74").Value = 820.2
$ws.Range("J74").Value = 700
$ws.Range("K74").Value = 820.2
$ws.Range("L74").Value = 700
$ws.Range("M74").Value = 53.79999999999995
$ws.Range("N74").Value = -2448
$ws.Range("H77").Value = 809.2727
$ws.Range("I77").Value = 820.2
$ws.Range("J77").Value = 700
$ws.Range("K77").Value = 4101
$ws.Range("L77").Value = 3500
$ws.Range("M77").Value = 267
$ws.Range("N77").Value = -12236
$ws.Range("H97").Value = 803.6667
$ws.Range("I97").Value = 624.5238
$ws.Range("J97").Value = 2057.6667
$ws.Range("K97").Value = 624.5238
$ws.Range("L97").Value = 2057.6667
$ws.Range("M97").Value = -128.5238000000001
$ws.Range("N97").Value = -3049.6667
$ws.Range("H122").Value = 3252.68
$ws.Range("I122").Value = 3092.6843
$ws.Range("J122").Value = 3759.3333
$ws.Range("K122").Value = 9278.052899999999
$ws.Range("L122").Value = 11277.9999
$ws.Range("M122").Value = -6828.052899999999
$ws.Range("N122").Value = -16177.9999
$ws.Range("H132").Value = 4961.15
$ws.Range("I132").Value = 6491.864
$ws.Range("J132").Value = 3090.2778
$ws.Range("K132").Value = 19475.592
$ws.Range("L132").Value = 9270.8334
$ws.Range("M132").Value = -16945.592
$ws.Range("N132").Value = -14330.8334
$ws.Range("H136").Value = 1831.6552
$ws.Range("I136").Value = 1751.9445
$ws.Range("J136").Value = 1962.091
$ws.Range("K136").Value = 5255.833500000001
$ws.Range("L136").Value = 5886.272999999999
$ws.Range("M136").Value = -2705.833500000001
$ws.Range("N136").Value = -10986.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 14619.25
$ws.Range("I22").Value = 18566.666
$ws.Range("K22").Value = 18566.666
$ws.Range("M22").Value = -18393.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 552.5
$ws.Range("I19").Value = 163
$ws.Range("K19").Value = 163
$ws.Range("M19").Value = 7
$ws.Range("H24").Value = 552.5
$ws.Range("I24").Value = 163
$ws.Range("K24").Value = 163
$ws.Range("M24").Value = 7
$ws.Range("H31").Value = 2363.3225
$ws.Range("I31").Value = 1210.25
$ws.Range("J31").Value = 4459.8184
$ws.Range("K31").Value = 1210.25
$ws.Range("L31").Value = 4459.8184
$ws.Range("M31").Value = -915.25
$ws.Range("N31").Value = -5049.8184
$ws.Range("H34").Value = 2363.3225
$ws.Range("I34").Value = 1210.25
$ws.Range("J34").Value = 4459.8184
$ws.Range("K34").Value = 1210.25
$ws.Range("L34").Value = 4459.8184
$ws.Range("M34").Value = -1008.25
$ws.Range("N34").Value = -4863.8184
$ws.Range("H122").Value = 2520.5264
$ws.Range("I122").Value = 2578.875
$ws.Range("K122").Value = 7736.625
$ws.Range("M122").Value = -5286.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 42884.582
$ws.Range("I121").Value = 879.8571
$ws.Range("J121").Value = 101691.2
$ws.Range("K121").Value = 2639.5713
$ws.Range("L121").Value = 305073.6
$ws.Range("M121").Value = -1329.5713
$ws.Range("N121").Value = -307693.6
$ws.Range("H123").Value = 6503.5
$ws.Range("I123").Value = 3007.5
$ws.Range("J123").Value = 9999.5
$ws.Range("K123").Value = 9022.5
$ws.Range("L123").Value = 29998.5
$ws.Range("M123").Value = -6572.5
$ws.Range("N123").Value = -34898.5
$ws.Range("H129").Value = 6250661
$ws.Range("I129").Value = 729.8
$ws.Range("J129").Value = 16667213
$ws.Range("K129").Value = 2189.4
$ws.Range("L129").Value = 50001639
$ws.Range("M129").Value = 2810.6
$ws.Range("N129").Value = -50011639
$ws.Range("H132").Value = 1785.174
$ws.Range("I132").Value = 1436.2354
$ws.Range("J132").Value = 2773.8333
$ws.Range("K132").Value = 12926.1186
$ws.Range("L132").Value = 24964.4997
$ws.Range("M132").Value = -10396.1186
$ws.Range("N132").Value = -30024.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 32500
$ws.Range("J62").Value = 32500
$ws.Range("L62").Value = 32500
$ws.Range("N62").Value = -33872
$ws.Range("H65").Value = 32500
$ws.Range("J65").Value = 32500
$ws.Range("L65").Value = 97500
$ws.Range("N65").Value = -104364
$ws.Range("H126").Value = 1800
$ws.Range("I126").Value = 750
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 2250
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = 220
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 2815.5
$ws.Range("I132").Value = 2409.5
$ws.Range("J132").Value = 3871.1
$ws.Range("K132").Value = 7228.5
$ws.Range("L132").Value = 11613.3
$ws.Range("M132").Value = -4698.5
$ws.Range("N132").Value = -16673.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3698.6
$ws.Range("I7").Value = 3498.1333
$ws.Range("J7").Value = 4300
$ws.Range("K7").Value = 3498.1333
$ws.Range("L7").Value = 4300
$ws.Range("M7").Value = -3386.1333
$ws.Range("N7").Value = -4524
$ws.Range("H46").Value = 1493.5625
$ws.Range("I46").Value = 1183.3334
$ws.Range("J46").Value = 1679.7
$ws.Range("K46").Value = 1183.3334
$ws.Range("L46").Value = 1679.7
$ws.Range("M46").Value = -995.3334
$ws.Range("N46").Value = -2055.7
$ws.Range("H100").Value = 2581.75
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2581.75
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2581.75
$ws.Range("N100").Value = -3663.75
$ws.Range("M100").ClearContents()
$ws.Range("H126").Value = 3698.6
$ws.Range("I126").Value = 3498.1333
$ws.Range("J126").Value = 4300
$ws.Range("K126").Value = 10494.3999
$ws.Range("L126").Value = 12900
$ws.Range("M126").Value = -8024.3999
$ws.Range("N126").Value = -17840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H96").Value = 2166.6667
$ws.Range("I96").Value = 2500
$ws.Range("K96").Value = 2500
$ws.Range("M96").Value = -1127
